$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update stock values in column C for rows 8, 10, 13, 14, 15, 16, 17 from 1 -> 5
$ws.Range("C8").Value = 5
$ws.Range("C10").Value = 5
$ws.Range("C13").Value = 5
$ws.Range("C14").Value = 5
$ws.Range("C15").Value = 5
$ws.Range("C16").Value = 5
$ws.Range("C17").Value = 5

# Update the active selection to C17
$ws.Range("C17").Select()
